$wb = $excel.ActiveWorkbook

# New "want to go" counts (column F) for rows 2-7
$updates = @{
    2 = 1361
    3 = 2014
    4 = 239
    5 = 69
    6 = 6356
    7 = 232
}

# Both "展览" and "全部类型" sheets contain identical data and both need updating
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
